# Update election result tables for PORTALEGRE / ELVAS (row 2) with the
# latest vote counts ("tabelas atualizadas, portugal, europa e fora da
# europa/ Graficos a funcionar").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 204
$ws.Range("I2").Value  = 550
$ws.Range("J2").Value  = 2315
$ws.Range("K2").Value  = 22
$ws.Range("L2").Value  = 591
$ws.Range("M2").Value  = 35
$ws.Range("N2").Value  = 388
$ws.Range("O2").Value  = 2
$ws.Range("P2").Value  = 10
$ws.Range("Q2").Value  = 3
$ws.Range("R2").Value  = 29
$ws.Range("S2").Value  = 264
$ws.Range("T2").Value  = 411
$ws.Range("U2").Value  = 42
$ws.Range("V2").Value  = 3485
$ws.Range("W2").Value  = 4
$ws.Range("X2").Value  = 3560
$ws.Range("Y2").Value  = 7
$ws.Range("Z2").Value  = 66
$ws.Range("AA2").Value = 23
